# Applies the "contingencies with rene fine" edit:
#  - Extends header row 1 with two new columns P1=14, Q1=15 (same style as existing header cells)
#  - Adds two new columns P and Q (value 2) for every data row 2..25
#  - Swaps values in existing columns I/K and M/O for every data row 2..25
#    (I: 1->2, K: 2->1, M: 1->2, O: 2->1)
#  - Dimension grows from A1:O25 to A1:Q25 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style (bold/border/center) used by the header row and first column -- same as cell O1.
$headerStyleRange = $ws.Range("O1")

# --- Header row additions: P1, Q1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$headerStyleRange.Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Per-row edits for rows 2 through 25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap I <-> K and M <-> O by setting their new values directly
    $ws.Cells.Item($r, 9).Value  = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1   # column K
    $ws.Cells.Item($r, 13).Value = 2   # column M
    $ws.Cells.Item($r, 15).Value = 1   # column O

    # Add new columns P (16) and Q (17) with value 2
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
